# Split the field-code paragraph's runs so that:
#   "{m"   ->  "{"  +  "m"
#   "()}"  ->  "()" +  "}"
#
# We locate the paragraph containing the M2Doc field text, then force a
# run boundary at each split point by temporarily inserting a paragraph
# break there and immediately deleting that paragraph mark again. This
# rejoins the two paragraphs back into one while leaving the text on
# either side of the break as two distinct runs (Word does not re-merge
# runs across a former paragraph-mark deletion the way it would merge
# runs touched by a plain Find/Replace or InsertBefore/InsertAfter).

$d = $word.ActiveDocument

function Split-RunAt($absPos) {
    $r = $d.Range($absPos, $absPos)
    $r.InsertParagraphAfter()
    $mark = $d.Range($absPos, $absPos + 1)
    $mark.Delete()
}

# Find the paragraph that contains the "{m" field opening.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.IndexOf("{m") -ge 0) {
        $target = $p
    }
}

$pStart = $target.Range.Start
$text = $target.Range.Text

# --- Split 1: "{m" -> "{" + "m" ---
$idx1 = $text.IndexOf("{m")
Split-RunAt ($pStart + $idx1 + 1)

# --- Split 2: "()}" -> "()" + "}" ---
# Re-read the paragraph text/start; splitting above did not change offsets
# before the split point, and "()}" is further along in the same paragraph.
$text2 = $target.Range.Text
$idx2 = $text2.IndexOf("()}")
Split-RunAt ($pStart + $idx2 + 2)
